# Apply updated result values (Update Name of Algo)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C2").Value = -12.463
$ws.Range("A3").Value = -21.674
$ws.Range("C5").Value = -12.836
$ws.Range("E5").Value = 13.034
$ws.Range("E9").Value = 13.378
$ws.Range("E11").Value = 13.068
$ws.Range("A14").Value = -20.814
$ws.Range("A16").Value = -21.215
$ws.Range("C16").Value = -12.09
$ws.Range("E17").Value = 13.784
$ws.Range("A21").Value = -21.04
$ws.Range("E21").Value = 13.535
$ws.Range("A23").Value = -21.709
$ws.Range("A25").Value = -22.269
